$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase practice trials from 2 to 40 (word_trial_count + nonword_trial_count)
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 20

# Update the selected cell on the sheet
$ws.Range("E3").Select()
